$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.924.20"
$ws.Range("E2").Value = "  +3.89%  "
$ws.Range("D3").Value = "2.469.05"
$ws.Range("E3").Value = "  +5.62%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.26"
$ws.Range("E5").Value = "  +2.57%  "
$ws.Range("E6").Value = "  +8.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  +1.25%  "
$ws.Range("D9").Value = "2.464.47"
$ws.Range("E10").Value = "  +2.63%  "
$ws.Range("E11").Value = "  +1.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.151"
$ws.Range("E12").Value = "  +1.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.57"
$ws.Range("E14").Value = "  +10.89%  "
$ws.Range("D15").Value = "2.909.60"
$ws.Range("E15").Value = "  +5.62%  "
$ws.Range("D16").Value = "62.797.70"
$ws.Range("E16").Value = "  +3.82%  "
$ws.Range("E17").Value = "  +4.47%  "
$ws.Range("D18").Value = "2.465.03"
$ws.Range("E18").Value = "  +5.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.23"
$ws.Range("E19").Value = "  +5.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "339.83"
$ws.Range("E20").Value = "  +7.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.26"
$ws.Range("E21").Value = "  +3.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.81"
$ws.Range("E22").Value = "  +3.03%  "
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.60"
$ws.Range("E24").Value = "  +2.06%  "
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("E27").Value = "  +6.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.05"
$ws.Range("E28").Value = "  +1.00%  "
$ws.Range("E29").Value = "  +9.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.84"
$ws.Range("E30").Value = "  +12.35%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0798"
$ws.Range("E31").Value = "  +8.49%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.84"
$ws.Range("E32").Value = "  +5.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "175.91"
$ws.Range("E33").Value = "  +2.79%  "
$ws.Range("E34").Value = "  +9.90%  "
$ws.Range("E35").Value = "  +2.55%  "
$ws.Range("E36").Value = "  +3.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "368.08"
$ws.Range("E37").Value = "  +11.66%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  +5.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.60"
$ws.Range("E42").Value = "  +6.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "149.58"
$ws.Range("E43").Value = "  +8.34%  "
$ws.Range("E44").Value = "  +4.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.43"
$ws.Range("E45").Value = "  +5.78%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.596"
$ws.Range("E46").Value = "  +4.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0958"
$ws.Range("E47").Value = "  +0.64%  "
$ws.Range("E48").Value = "  +2.98%  "
$ws.Range("D49").Value = "0.0₆0232"
$ws.Range("E49").Value = "  +3.93%  "
$ws.Range("E50").Value = "  +4.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.91"
$ws.Range("E51").Value = "  +4.25%  "
